# Apply the changes described in the diff:
# - Update several "Valor PAA" (column D) figures for year 2022 rows
# - Move the active cell selection from F2 to D25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D14").Value = 17886815132
$ws.Range("D15").Value = 17588104219
$ws.Range("D16").Value = 1274763523
$ws.Range("D17").Value = 7492465876
$ws.Range("D18").Value = 24479997671
$ws.Range("D20").Value = 19925748000
$ws.Range("D21").Value = 140046000
$ws.Range("D22").Value = 1708649347
$ws.Range("D23").Value = 2104

$ws.Range("D25").Select()
